$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-06 Sunday", "2024-10-07 Monday"),
    @("55×67=3685", "87×56=4872"),
    @("37×22=814", "60×16=960"),
    @("91×92=8372", "20×11=220"),
    @("97×48=4656", "21×76=1596"),
    @("32×13=416", "39×70=2730"),
    @("96×89=8544", "15×36=540"),
    @("81×96=7776", "83×47=3901"),
    @("87×52=4524", "30×71=2130"),
    @("80×66=5280", "56×96=5376"),
    @("32×78=2496", "15×41=615"),
    @("74×22=1628", "85×97=8245"),
    @("78×31=2418", "17×34=578"),
    @("33×62=2046", "24×24=576"),
    @("69×25=1725", "22×27=594"),
    @("19×47=893", "82×92=7544"),
    @("96×12=1152", "63×49=3087"),
    @("28×22=616", "75×26=1950"),
    @("25×65=1625", "67×67=4489"),
    @("84×15=1260", "45×19=855"),
    @("99×64=6336", "19×81=1539"),
    @("23×16=368", "79×98=7742"),
    @("79×43=3397", "22×19=418"),
    @("67×44=2948", "93×45=4185"),
    @("28×76=2128", "59×45=2655"),
    @("76×77=5852", "68×72=4896")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
